$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.291.64'
$ws.Range('E2').Value = '  +2.00%  '
$ws.Range('D3').Value = '1.845.56'
$ws.Range('E3').Value = '  +2.02%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.29'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.612'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.36%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '43.20'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +15.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.306'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.70%  '
$ws.Range('E10').Value = '  +1.80%  '
$ws.Range('E11').Value = '  +3.47%  '
$ws.Range('D12').Value = '2.110.66'
$ws.Range('E12').Value = '  +1.94%  '
$ws.Range('E13').Value = '  +2.30%  '
$ws.Range('D14').Value = '1.836.57'
$ws.Range('E14').Value = '  +1.51%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.74'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.88%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.660'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.03%  '
$ws.Range('D17').Value = '35.235.04'
$ws.Range('E17').Value = '  +2.01%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.77'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.54%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '246.09'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.05%  '
$ws.Range('D20').Value = '0.0₃0794'
$ws.Range('E20').Value = '  +2.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.11'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +7.90%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.71'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +14.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.18'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.61%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '172.15'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.91'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.92'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.48%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.123'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.43%  '
$ws.Range('D29').Value = '3.646.17'
$ws.Range('E29').Value = '  +50.07%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.19%  '
$ws.Range('E31').Value = '  +8.38%  '
$ws.Range('E32').Value = '  +3.50%  '
$ws.Range('E33').Value = '  +2.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0538'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.79%  '
$ws.Range('E35').Value = '  +4.31%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.674'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '90.43'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +11.73%  '
$ws.Range('E38').Value = '  +0.95%  '
$ws.Range('D39').Value = '1.339.94'
$ws.Range('E39').Value = '  -1.91%  '
$ws.Range('E40').Value = '  +8.73%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.42'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.94%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0194'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '14.83'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +7.32%  '
$ws.Range('E44').Value = '  +5.96%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.84'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.76%  '
$ws.Range('E46').Value = '  +0.76%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0519'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.50%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.07'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.17%  '
$ws.Range('D49').Value = '2.011.02'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '104.45'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.04%  '
